$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.056.31"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "1.816.35"
$ws.Range("E3").Value = "  +1.36%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.74"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4993"
$ws.Range("E7").Value = "  -2.58%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3908"
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09851"
$ws.Range("E9").Value = "  +26.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.108"
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.426"
$ws.Range("E12").Value = "  +2.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.57"
$ws.Range("E13").Value = "  +1.75%  "
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").Value = "1.811.02"
$ws.Range("E15").Value = "  +1.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.274"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("E17").Value = "  +5.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.35"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06641"
$ws.Range("E19").Value = "  +1.69%  "
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.20"
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.953"
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("D23").Value = "28.103.98"
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.25"
$ws.Range("E24").Value = "  +1.98%  "
$ws.Range("E25").Value = "  +0.74%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.10"
$ws.Range("E26").Value = "  -0.96%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.66"
$ws.Range("E27").Value = "  +1.60%  "
$ws.Range("B28").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C28").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D28").Value = "2.024.50"
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.407"
$ws.Range("E29").Value = "  +1.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.68"
$ws.Range("E30").Value = "  -0.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1062"
$ws.Range("E31").Value = "  -1.65%  "
$ws.Range("E32").Value = "  -0.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.566"
$ws.Range("E33").Value = "  +1.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.617"
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("E35").Value = "  -5.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02338"
$ws.Range("E36").Value = "  +1.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.900"
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2141"
$ws.Range("E38").Value = "  +0.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.956"
$ws.Range("E39").Value = "  -1.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.36"
$ws.Range("E40").Value = "  -1.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6193"
$ws.Range("E41").Value = "  +1.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.183"
$ws.Range("E42").Value = "  +2.66%  "
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.13"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5900"
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.692"
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("B47").Value = "WEMIXTOKEN"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.274"
$ws.Range("E47").Value = "  -2.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.37"
$ws.Range("E48").Value = "  -0.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.940"
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("E50").Value = "  -2.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06782"
$ws.Range("E51").Value = "  -0.51%  "
